$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the D column (Price) range to Text format first so numeric-looking
# strings (e.g. "1.001", "324.60") are preserved exactly as strings rather
# than being coerced into numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.917.90"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.907.65"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "324.60"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.4589"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").Value = "0.3819"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "0.9794"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "22.20"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("D12").Value = "1.883.56"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").Value = "6.957"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "5.678"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").Value = "0.07072"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "83.85"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "0.000009457"
$ws.Range("D19").Value = "16.65"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "28.915.93"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "2.096"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "158.66"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "19.06"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "5.674"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "117.64"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("D29").Value = "1.876"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "0.09306"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").Value = "0.8629"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D33").Value = "1.246"
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("D34").Value = "3.027"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "0.05715"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "1.158"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "0.02045"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "7.471"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "0.5497"
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("D43").Value = "9.346"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "2.173"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").Value = "0.000002748"
$ws.Range("E45").Value = "  -11.34%  "
$ws.Range("D46").Value = "0.5183"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "11.23"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "0.06885"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").Value = "1.780"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "110.57"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "0.2861"
$ws.Range("E51").Value = "  -3.72%  "

# MXToken (row 41) and Algorand (row 42) swapped positions in the ranking
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.1755"
$ws.Range("E41").Value = "  -1.44%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.915"
$ws.Range("E42").Value = "  +7.70%  "

# Clean up the temporary Text number format so the cell style index
# matches the original workbook (no explicit style / quote-prefix left behind).
$priceRange.Style = "Normal"
